$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 690; this shifts the existing rows 690-730 down to 691-731
# and grows the sheet's used range from A1:R730 to A1:R731 (matching the diff).
$ws.Rows.Item(690).Insert()

# Populate the newly inserted row 690 with the new daily price record.
$ws.Range("A690").Value = 3
$ws.Range("B690").Value = "Femacal de La Calera"
$ws.Range("C690").Value = "Coquimbo"
$ws.Range("D690").Value = 45021
$ws.Range("E690").Value = 5
$ws.Range("F690").Value = 100112003
$ws.Range("G690").Value = "Ajo"
$ws.Range("H690").Value = "Chino"
$ws.Range("I690").Value = "Primera"
$ws.Range("J690").Value = 73
$ws.Range("K690").Value = 14500
$ws.Range("L690").Value = 15000
$ws.Range("M690").Value = 14760
$ws.Range("N690").Value = "$/caja 10 kilos"
$ws.Range("O690").Value = "China"
$ws.Range("P690").Value = 1476
$ws.Range("Q690").Value = 10
$ws.Range("R690").Value = "Hortaliza"
